# Generate Report for Handback
# Updates the localization-status workbook to reflect a handback event:
#  - Status cells move from "In Translation" to "Handed back: in sync with en-US"
#  - Latest Handback DateTime is stamped for zh-cn and de-de
#  - Latest Target File / Latest Handback File columns (I/J) are populated
#    with hyperlinks + filenames for each handed-back source file
#  - A few columns are widened so the new content is readable

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet: zh-cn/de-de status columns (E, F) for all 3 rows ----
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Range("E4").Value = $statusText
$overview.Range("F4").Value = $statusText

# ---- zh-cn / de-de sheets: Status column (C) for all 3 rows ----
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText
$zhcn.Range("C4").Value = $statusText

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText
$dede.Range("C4").Value = $statusText

# ---- zh-cn sheet: Latest Handback DateTime (K) ----
$zhcn.Range("K2").Value = "2016-09-06 17:09:20"
$zhcn.Range("K3").Value = "2016-09-06 17:09:20"
$zhcn.Range("K4").Value = "2016-09-06 17:09:20"

# ---- de-de sheet: Latest Handback DateTime (K) ----
$dede.Range("K2").Value = "2016-09-06 17:09:40"
$dede.Range("K3").Value = "2016-09-06 17:09:40"
$dede.Range("K4").Value = "2016-09-06 17:09:40"

# ---- zh-cn sheet: Latest Target File (I) + Latest Handback File (J) ----
$zhcn.Range("I2").Value = "61985fc3-8070-4228-ae4c-7c06387bbbef.yml"
$zhcn.Range("J2").Value = "61985fc3-8070-4228-ae4c-7c06387bbbef.f936b3f3ec1df9529148f734b1d064c589fefce3.zh-cn.xlf"

$zhcn.Range("I3").Value = "85e3c3a8-f6f5-4df3-b00b-6d52a2f4fdf2.yml"
$zhcn.Range("J3").Value = "85e3c3a8-f6f5-4df3-b00b-6d52a2f4fdf2.834c46732b496485a33f3ef5f31abcf6be6667eb.zh-cn.xlf"

$zhcn.Range("I4").Value = "acc813c1-5faa-4f55-8143-732558871e84.md"
$zhcn.Range("J4").Value = "acc813c1-5faa-4f55-8143-732558871e84.7a6f9709385fec5a8d82071ac6e1aaeade9ef8c3.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e3f8705a5bc5196548c9473e6a86560464b476e/e2e/61985fc3-8070-4228-ae4c-7c06387bbbef.yml", "", "", "61985fc3-8070-4228-ae4c-7c06387bbbef.yml")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e3f8705a5bc5196548c9473e6a86560464b476e/e2e/85e3c3a8-f6f5-4df3-b00b-6d52a2f4fdf2.yml", "", "", "85e3c3a8-f6f5-4df3-b00b-6d52a2f4fdf2.yml")
$zhcn.Hyperlinks.Add($zhcn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e3f8705a5bc5196548c9473e6a86560464b476e/e2e/acc813c1-5faa-4f55-8143-732558871e84.md", "", "", "acc813c1-5faa-4f55-8143-732558871e84.md")

# ---- de-de sheet: Latest Target File (I) + Latest Handback File (J) ----
$dede.Range("I2").Value = "61985fc3-8070-4228-ae4c-7c06387bbbef.yml"
$dede.Range("J2").Value = "61985fc3-8070-4228-ae4c-7c06387bbbef.f936b3f3ec1df9529148f734b1d064c589fefce3.de-de.xlf"

$dede.Range("I3").Value = "85e3c3a8-f6f5-4df3-b00b-6d52a2f4fdf2.yml"
$dede.Range("J3").Value = "85e3c3a8-f6f5-4df3-b00b-6d52a2f4fdf2.834c46732b496485a33f3ef5f31abcf6be6667eb.de-de.xlf"

$dede.Range("I4").Value = "acc813c1-5faa-4f55-8143-732558871e84.md"
$dede.Range("J4").Value = "acc813c1-5faa-4f55-8143-732558871e84.7a6f9709385fec5a8d82071ac6e1aaeade9ef8c3.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e3f8705a5bc5196548c9473e6a86560464b476e/e2e/61985fc3-8070-4228-ae4c-7c06387bbbef.yml", "", "", "61985fc3-8070-4228-ae4c-7c06387bbbef.yml")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e3f8705a5bc5196548c9473e6a86560464b476e/e2e/85e3c3a8-f6f5-4df3-b00b-6d52a2f4fdf2.yml", "", "", "85e3c3a8-f6f5-4df3-b00b-6d52a2f4fdf2.yml")
$dede.Hyperlinks.Add($dede.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e3f8705a5bc5196548c9473e6a86560464b476e/e2e/acc813c1-5faa-4f55-8143-732558871e84.md", "", "", "acc813c1-5faa-4f55-8143-732558871e84.md")

# ---- Column width adjustments (widen to fit the newly populated content) ----
$overview.Range("E1").ColumnWidth = 29.1
$overview.Range("F1").ColumnWidth = 29.1

$zhcn.Range("C1").ColumnWidth = 29.1
$zhcn.Range("I1").ColumnWidth = 39.15
$zhcn.Range("J1").ColumnWidth = 39.15

$dede.Range("C1").ColumnWidth = 29.1
$dede.Range("I1").ColumnWidth = 39.15
$dede.Range("J1").ColumnWidth = 39.15
